# Applies the scheduled-runner update to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4411.1113
$ws.Range("I64").Value = 5400
$ws.Range("J64").Value = 3175
$ws.Range("K64").Value = 5400
$ws.Range("L64").Value = 3175
$ws.Range("M64").Value = -5152
$ws.Range("N64").Value = -3671
$ws.Range("H67").Value = 4411.1113
$ws.Range("I67").Value = 5400
$ws.Range("J67").Value = 3175
$ws.Range("K67").Value = 5400
$ws.Range("L67").Value = 3175
$ws.Range("M67").Value = -4542
$ws.Range("N67").Value = -4891
$ws.Range("H132").Value = 2916.111
$ws.Range("I132").Value = 2916.111
$ws.Range("K132").Value = 8748.332999999999
$ws.Range("M132").Value = -6218.332999999999
$ws.Range("H137").Value = 2567.4333
$ws.Range("I137").Value = 1286.2941
$ws.Range("K137").Value = 3858.8823
$ws.Range("M137").Value = -1308.8823
$ws.Range("H141").Value = 3709.9285
$ws.Range("I141").Value = 3072.2307
$ws.Range("K141").Value = 9216.6921
$ws.Range("M141").Value = -4036.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10523.212
$ws.Range("I32").Value = 7938.107
$ws.Range("J32").Value = 24999.8
$ws.Range("K32").Value = 7938.107
$ws.Range("L32").Value = 24999.8
$ws.Range("M32").Value = -7651.107
$ws.Range("N32").Value = -25573.8
$ws.Range("H45").Value = 1850.3334
$ws.Range("I45").Value = 1768.5
$ws.Range("K45").Value = 1768.5
$ws.Range("M45").Value = -1391.5
$ws.Range("H61").Value = 4894.2
$ws.Range("I61").Value = 4868.25
$ws.Range("K61").Value = 4868.25
$ws.Range("M61").Value = -4656.25
$ws.Range("H74").Value = 1752.6394
$ws.Range("I74").Value = 1227.86
$ws.Range("K74").Value = 1227.86
$ws.Range("M74").Value = -353.8599999999999
$ws.Range("H77").Value = 1752.6394
$ws.Range("I77").Value = 1227.86
$ws.Range("K77").Value = 6139.299999999999
$ws.Range("M77").Value = -1771.299999999999
$ws.Range("H122").Value = 2737.9
$ws.Range("I122").Value = 1893.4286
$ws.Range("K122").Value = 5680.2858
$ws.Range("M122").Value = -3230.2858
$ws.Range("H132").Value = 2150.8
$ws.Range("J132").Value = 3999
$ws.Range("L132").Value = 11997
$ws.Range("N132").Value = -17057

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 99999
$ws.Range("J82").Value = 99999
$ws.Range("L82").Value = 99999
$ws.Range("N82").Value = -100765
$ws.Range("H85").Value = 99999
$ws.Range("J85").Value = 99999
$ws.Range("L85").Value = 99999
$ws.Range("N85").Value = -102651
$ws.Range("H86").Value = 1742.3
$ws.Range("I86").Value = 1807.1666
$ws.Range("K86").Value = 1807.1666
$ws.Range("M86").Value = -684.1666
$ws.Range("H89").Value = 1742.3
$ws.Range("I89").Value = 1807.1666
$ws.Range("K89").Value = 9035.833000000001
$ws.Range("M89").Value = -3419.833000000001
$ws.Range("H134").Value = 2139.7646
$ws.Range("I134").Value = 1628.1538
$ws.Range("J134").Value = 3802.5
$ws.Range("K134").Value = 4884.4614
$ws.Range("L134").Value = 11407.5
$ws.Range("M134").Value = -2349.4614
$ws.Range("N134").Value = -16477.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 376.4
$ws.Range("I22").Value = 177.16667
$ws.Range("K22").Value = 177.16667
$ws.Range("M22").Value = 172.83333
$ws.Range("H31").Value = 4218.4443
$ws.Range("I31").Value = 2186.6
$ws.Range("J31").Value = 6758.25
$ws.Range("K31").Value = 2186.6
$ws.Range("L31").Value = 6758.25
$ws.Range("M31").Value = -1891.6
$ws.Range("N31").Value = -7348.25
$ws.Range("H34").Value = 4218.4443
$ws.Range("I34").Value = 2186.6
$ws.Range("J34").Value = 6758.25
$ws.Range("K34").Value = 2186.6
$ws.Range("L34").Value = 6758.25
$ws.Range("M34").Value = -1984.6
$ws.Range("N34").Value = -7162.25
$ws.Range("H64").Value = 19749.75
$ws.Range("J64").Value = 19749.75
$ws.Range("L64").Value = 19749.75
$ws.Range("N64").Value = -20245.75
$ws.Range("H67").Value = 19749.75
$ws.Range("J67").Value = 19749.75
$ws.Range("L67").Value = 19749.75
$ws.Range("N67").Value = -21465.75
$ws.Range("H134").Value = 2712.9688
$ws.Range("I134").Value = 1930.4286
$ws.Range("J134").Value = 4206.909
$ws.Range("K134").Value = 5791.2858
$ws.Range("L134").Value = 12620.727
$ws.Range("M134").Value = -3256.2858
$ws.Range("N134").Value = -17690.727
$ws.Range("H138").Value = 100000
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 671.4286
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 900
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 2700
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -5196
$ws.Range("H130").Value = 3000
$ws.Range("J130").Value = 4000
$ws.Range("L130").Value = 12000
$ws.Range("N130").Value = -22040
$ws.Range("H131").Value = 1457.3914
$ws.Range("I131").Value = 1259.75
$ws.Range("K131").Value = 3779.25
$ws.Range("M131").Value = 1260.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 79944.62
$ws.Range("I122").Value = 2408.6
$ws.Range("J122").Value = 338398
$ws.Range("K122").Value = 7225.799999999999
$ws.Range("L122").Value = 1015194
$ws.Range("M122").Value = -4775.799999999999
$ws.Range("N122").Value = -1020094

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 110
$ws.Range("I22").Value = 100
$ws.Range("J22").Value = 120
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = 120
$ws.Range("M22").Value = 195
$ws.Range("N22").Value = -710
$ws.Range("H27").Value = 110
$ws.Range("I27").Value = 100
$ws.Range("J27").Value = 120
$ws.Range("K27").Value = 100
$ws.Range("L27").Value = 120
$ws.Range("M27").Value = 7
$ws.Range("N27").Value = -334
$ws.Range("H46").Value = 1993.1111
$ws.Range("I46").Value = 1756.3334
$ws.Range("K46").Value = 1756.3334
$ws.Range("M46").Value = -1568.3334
$ws.Range("H55").Value = 331.29413
$ws.Range("I55").Value = 340.92856
$ws.Range("K55").Value = 340.92856
$ws.Range("M55").Value = -167.92856
$ws.Range("H68").Value = 6944.3335
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 6944.3335
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 6944.3335
$ws.Range("M68").Value = ""
$ws.Range("N68").Value = -8442.333500000001
$ws.Range("H71").Value = 6944.3335
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 6944.3335
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 34721.6675
$ws.Range("M71").Value = ""
$ws.Range("N71").Value = -42209.6675
$ws.Range("H93").Value = 2673.7856
$ws.Range("I93").Value = 2673.7856
$ws.Range("K93").Value = 2673.7856
$ws.Range("M93").Value = -1425.7856
$ws.Range("H122").Value = 3831.6667
$ws.Range("I122").Value = 3841.875
$ws.Range("K122").Value = 11525.625
$ws.Range("M122").Value = -9075.625
$ws.Range("H133").Value = 20000
$ws.Range("J133").Value = 20000
$ws.Range("L133").Value = 20000
$ws.Range("N133").Value = -25060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 535.3333
$ws.Range("I107").Value = 377.75
$ws.Range("J107").Value = 592.63635
$ws.Range("K107").Value = 1133.25
$ws.Range("L107").Value = 1777.90905
$ws.Range("M107").Value = 786.75
$ws.Range("N107").Value = -5617.90905

